# Auto-generated edit script: update Typhon_Profits market price data
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5436504.5
$ws.Range("J40").Value = 6212862.5
$ws.Range("L40").Value = 6212862.5
$ws.Range("N40").Value = -6213212.5
$ws.Range("H98").Value = 1773
$ws.Range("I98").Value = 1163.3334
$ws.Range("J98").Value = 2687.5
$ws.Range("K98").Value = 1163.3334
$ws.Range("L98").Value = 2687.5
$ws.Range("M98").Value = 334.6666
$ws.Range("N98").Value = -5683.5
$ws.Range("H103").Value = 250401.2
$ws.Range("I103").Value = 500082
$ws.Range("K103").Value = 1500246
$ws.Range("M103").Value = -1499660
$ws.Range("H122").Value = 1773
$ws.Range("I122").Value = 1163.3334
$ws.Range("J122").Value = 2687.5
$ws.Range("K122").Value = 3490.0002
$ws.Range("L122").Value = 8062.5
$ws.Range("M122").Value = -1040.0002
$ws.Range("N122").Value = -12962.5
$ws.Range("H132").Value = 3515.125
$ws.Range("I132").Value = 3888.7144
$ws.Range("K132").Value = 11666.1432
$ws.Range("M132").Value = -9136.143199999999
$ws.Range("H137").Value = 2455.4285
$ws.Range("I137").Value = 2211.8572
$ws.Range("K137").Value = 6635.571599999999
$ws.Range("M137").Value = -4085.571599999999
$ws.Range("H138").Value = 2094.2222
$ws.Range("J138").Value = 2197.5874
$ws.Range("L138").Value = 6592.762199999999
$ws.Range("N138").Value = -16872.7622

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4764.4517
$ws.Range("I32").Value = 4874.6226
$ws.Range("K32").Value = 4874.6226
$ws.Range("M32").Value = -4587.6226
$ws.Range("H61").Value = 2244.75
$ws.Range("J61").Value = 3616.5
$ws.Range("L61").Value = 3616.5
$ws.Range("N61").Value = -4040.5
$ws.Range("H97").Value = 1162.3226
$ws.Range("I97").Value = 1210.7916
$ws.Range("J97").Value = 996.1429000000001
$ws.Range("K97").Value = 1210.7916
$ws.Range("L97").Value = 996.1429000000001
$ws.Range("M97").Value = -714.7916
$ws.Range("N97").Value = -1988.1429
$ws.Range("H132").Value = 31209.568
$ws.Range("I132").Value = 1678.5555
$ws.Range("J132").Value = 102084
$ws.Range("K132").Value = 5035.666499999999
$ws.Range("L132").Value = 306252
$ws.Range("M132").Value = -2505.666499999999
$ws.Range("N132").Value = -311312
$ws.Range("H136").Value = 2244.75
$ws.Range("J136").Value = 3616.5
$ws.Range("L136").Value = 10849.5
$ws.Range("N136").Value = -15949.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1859.4445
$ws.Range("I86").Value = 1746.0416
$ws.Range("K86").Value = 1746.0416
$ws.Range("M86").Value = -623.0416
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H89").Value = 1859.4445
$ws.Range("I89").Value = 1746.0416
$ws.Range("K89").Value = 8730.208000000001
$ws.Range("M89").Value = -3114.208000000001
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H105").Value = 4366.7
$ws.Range("I105").Value = 5533.6
$ws.Range("K105").Value = 5533.6
$ws.Range("M105").Value = -3786.6
$ws.Range("H134").Value = 6229.4736
$ws.Range("I134").Value = 6824.1333
$ws.Range("J134").Value = 3999.5
$ws.Range("K134").Value = 20472.3999
$ws.Range("L134").Value = 11998.5
$ws.Range("M134").Value = -17937.3999
$ws.Range("N134").Value = -17068.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14467.182
$ws.Range("I31").Value = 22236.947
$ws.Range("J31").Value = 3922.5
$ws.Range("K31").Value = 22236.947
$ws.Range("L31").Value = 3922.5
$ws.Range("M31").Value = -21941.947
$ws.Range("N31").Value = -4512.5
$ws.Range("H34").Value = 14467.182
$ws.Range("I34").Value = 22236.947
$ws.Range("J34").Value = 3922.5
$ws.Range("K34").Value = 22236.947
$ws.Range("L34").Value = 3922.5
$ws.Range("M34").Value = -22034.947
$ws.Range("N34").Value = -4326.5
$ws.Range("H58").Value = 25231.523
$ws.Range("I58").Value = 1529.9333
$ws.Range("J58").Value = 84485.5
$ws.Range("K58").Value = 1529.9333
$ws.Range("L58").Value = 84485.5
$ws.Range("M58").Value = -1326.9333
$ws.Range("N58").Value = -84891.5
$ws.Range("H132").Value = 26390.5
$ws.Range("I132").Value = 30093.277
$ws.Range("J132").Value = 9728
$ws.Range("K132").Value = 90279.83099999999
$ws.Range("L132").Value = 29184
$ws.Range("M132").Value = -87749.83099999999
$ws.Range("N132").Value = -34244
$ws.Range("H134").Value = 693.90625
$ws.Range("I134").Value = 578.75
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 1736.25
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = 798.75
$ws.Range("N134").Value = -9570
$ws.Range("H136").Value = 25231.523
$ws.Range("I136").Value = 1529.9333
$ws.Range("J136").Value = 84485.5
$ws.Range("K136").Value = 4589.7999
$ws.Range("L136").Value = 253456.5
$ws.Range("M136").Value = -2039.7999
$ws.Range("N136").Value = -258556.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 176.3
$ws.Range("I18").Value = 154.125
$ws.Range("J18").Value = 265
$ws.Range("K18").Value = 462.375
$ws.Range("L18").Value = 795
$ws.Range("M18").Value = -293.375
$ws.Range("N18").Value = -1133
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("M108").ClearContents()
$ws.Range("H122").Value = 682.5599999999999
$ws.Range("I122").Value = 360.875
$ws.Range("J122").Value = 833.94116
$ws.Range("K122").Value = 3247.875
$ws.Range("L122").Value = 7505.47044
$ws.Range("M122").Value = -797.875
$ws.Range("N122").Value = -12405.47044
$ws.Range("H128").Value = 310000
$ws.Range("I128").Value = 310000
$ws.Range("K128").Value = 930000
$ws.Range("M128").Value = -925020
$ws.Range("H131").Value = 810.88
$ws.Range("I131").Value = 799.5
$ws.Range("J131").Value = 811.11224
$ws.Range("K131").Value = 2398.5
$ws.Range("L131").Value = 2433.33672
$ws.Range("M131").Value = 2641.5
$ws.Range("N131").Value = -12513.33672

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 757.55554
$ws.Range("I97").Value = 742.5625
$ws.Range("J97").Value = 877.5
$ws.Range("K97").Value = 742.5625
$ws.Range("L97").Value = 877.5
$ws.Range("M97").Value = -246.5625
$ws.Range("N97").Value = -1869.5
$ws.Range("H122").Value = 333334080
$ws.Range("I122").Value = 111112100
$ws.Range("K122").Value = 333336300
$ws.Range("M122").Value = -333333850
$ws.Range("H132").Value = 16056.368
$ws.Range("I132").Value = 3183.3076
$ws.Range("K132").Value = 9549.9228
$ws.Range("M132").Value = -7019.9228

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 37703.57
$ws.Range("I136").Value = 56971.668
$ws.Range("K136").Value = 170915.004
$ws.Range("M136").Value = -168365.004

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 83335000
$ws.Range("I81").Value = 1615.8889
$ws.Range("J81").Value = 333335170
$ws.Range("K81").Value = 3231.7778
$ws.Range("L81").Value = 666670340
$ws.Range("M81").Value = -2170.7778
$ws.Range("N81").Value = -666672462
$ws.Range("H84").Value = 83335000
$ws.Range("I84").Value = 1615.8889
$ws.Range("J84").Value = 333335170
$ws.Range("K84").Value = 16158.889
$ws.Range("L84").Value = 3333351700
$ws.Range("M84").Value = -10854.889
$ws.Range("N84").Value = -3333362308
$ws.Range("H100").Value = 451.5
$ws.Range("I100").Value = 446.1111
$ws.Range("K100").Value = 892.2222
$ws.Range("M100").Value = -351.2222
$ws.Range("H122").Value = 1385.9
$ws.Range("J122").Value = 1549.8334
$ws.Range("L122").Value = 4649.5002
$ws.Range("N122").Value = -9549.5002
